# RECO_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer from 2021-05-12 to 2021-05-13
#  - refresh the Weight (D) / Percent Change (E) figures for rows 2-38
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected, so it must be unprotected before any cell can be written.
$ws.Unprotect()

$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.03238032964987903
$ws.Range("E2").Value = 0.01426583159359773
$ws.Range("D3").Value = 0.02849394166494817
$ws.Range("E3").Value = 0.004384757143764961
$ws.Range("D4").Value = 0.02853196658233451
$ws.Range("E4").Value = -0.02076634512325837
$ws.Range("D5").Value = 0.06341389317814874
$ws.Range("E5").Value = 0.003023534711955156
$ws.Range("D6").Value = 0.016052550838208
$ws.Range("E6").Value = 0.005790344412693882
$ws.Range("D7").Value = 0.01504338160217688
$ws.Range("E7").Value = 0.01998074145402007
$ws.Range("D8").Value = 0.03027869850163718
$ws.Range("E8").Value = 0.004943587290196616
$ws.Range("D9").Value = 0.03392064058908444
$ws.Range("E9").Value = 0.0131672597864767
$ws.Range("D10").Value = 0.02927194354607503
$ws.Range("E10").Value = 0.003464060373623568
$ws.Range("D11").Value = 0.03096918271576374
$ws.Range("E11").Value = -0.0628467485220553
$ws.Range("D12").Value = 0.0110523747969238
$ws.Range("E12").Value = -0.01483571493583347
$ws.Range("D13").Value = 0.01394810302941902
$ws.Range("E13").Value = 0.02261712439418417
$ws.Range("D14").Value = 0.0144086270288758
$ws.Range("E14").Value = -0.02121004789365666
$ws.Range("D15").Value = 0.009042285116463884
$ws.Range("E15").Value = 0.02109291562833748
$ws.Range("D16").Value = 0.007990664782186635
$ws.Range("E16").Value = 0.0258327668252889
$ws.Range("D17").Value = 0.02965098557970395
$ws.Range("E17").Value = 0.006269592476489061
$ws.Range("D18").Value = 0.02584889622114809
$ws.Range("E18").Value = -0.000684931506849229
$ws.Range("D19").Value = 0.03358867702460051
$ws.Range("E19").Value = -0.002336028751123043
$ws.Range("D20").Value = 0.03043502316200325
$ws.Range("E20").Value = 0.008957197157494656
$ws.Range("D21").Value = 0.04504805927058551
$ws.Range("E21").Value = 0.01022294871107765
$ws.Range("D22").Value = 0.03565087492518246
$ws.Range("E22").Value = 0.01244356659142221
$ws.Range("D23").Value = 0.03188862119436467
$ws.Range("E23").Value = 0.02656151419558372
$ws.Range("D24").Value = 0.03135788187126855
$ws.Range("E24").Value = 0.01622589213535042
$ws.Range("D25").Value = 0.01413259429525644
$ws.Range("E25").Value = 0.005253042921204543
$ws.Range("D26").Value = 0.0145955325751822
$ws.Range("E26").Value = 0.01248862790505334
$ws.Range("D27").Value = 0.0312709677743855
$ws.Range("E27").Value = 0.02031139419674433
$ws.Range("D28").Value = 0.03171881680138016
$ws.Range("E28").Value = -0.006647384178210758
$ws.Range("D29").Value = 0.02885065160423907
$ws.Range("E29").Value = 0.01686192468619252
$ws.Range("D30").Value = 0.02952041324434027
$ws.Range("E30").Value = 0.001649299047904673
$ws.Range("D31").Value = 0.03321687783237852
$ws.Range("E31").Value = -0.006777628375186251
$ws.Range("D32").Value = 0.03191940327034409
$ws.Range("E32").Value = 0.01483741246619985
$ws.Range("D33").Value = 0.02830200636766474
$ws.Range("E33").Value = 0.004990296645411618
$ws.Range("D34").Value = 0.03262256245693276
$ws.Range("E34").Value = 0.007869353923575906
$ws.Range("D35").Value = 0.03123837498805434
$ws.Range("E35").Value = 0.0002318571759796217
$ws.Range("D36").Value = 0.03167737165332944
$ws.Range("E36").Value = 0.02572245157192765
$ws.Range("D37").Value = 0.03266682426553062
$ws.Range("E37").Value = 0.02783799763500205
$ws.Range("D38").Value = 0.9999999999999999
$ws.Range("E38").Value = 0.005790249324755914
